$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.879.72"
$ws.Range("E2").Value = "  +1.70%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.891.83"
$ws.Range("E3").Value = "  +1.59%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.62%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.41"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4580"
$ws.Range("E7").Value = "  +0.55%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3903"
$ws.Range("E8").Value = "  +1.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07843"
$ws.Range("E9").Value = "  +0.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9886"
$ws.Range("E10").Value = "  +0.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.86"
$ws.Range("E11").Value = "  +1.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.932.77"
$ws.Range("E12").Value = "  +4.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.029"
$ws.Range("E13").Value = "  +2.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.708"
$ws.Range("E14").Value = "  +1.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06950"
$ws.Range("E15").Value = "  +0.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.94"
$ws.Range("E16").Value = "  +1.42%  "
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009970"
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.90"
$ws.Range("E19").Value = "  +1.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.901.74"
$ws.Range("E21").Value = "  +1.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.292"
$ws.Range("E22").Value = "  +1.01%  "
$ws.Range("E23").Value = "  +1.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.109.30"
$ws.Range("E24").Value = "  +1.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.060"
$ws.Range("E25").Value = "  -1.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.00"
$ws.Range("E26").Value = "  +1.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.30"
$ws.Range("E27").Value = "  +1.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.940"
$ws.Range("E28").Value = "  +5.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.930"
$ws.Range("E29").Value = "  +0.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.60"
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09355"
$ws.Range("E31").Value = "  +0.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9107"
$ws.Range("E32").Value = "  +0.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.298"
$ws.Range("E33").Value = "  +0.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.331"
$ws.Range("E34").Value = "  +1.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.267"
$ws.Range("E35").Value = "  -0.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.189"
$ws.Range("E36").Value = "  +3.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05767"
$ws.Range("E37").Value = "  +1.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02073"
$ws.Range("E38").Value = "  +0.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.000"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.740"
$ws.Range("E40").Value = "  +1.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5684"
$ws.Range("E41").Value = "  +2.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1772"
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.766"
$ws.Range("E43").Value = "  +1.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.264"
$ws.Range("E44").Value = "  +6.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "11.95"
$ws.Range("E45").Value = "  +3.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5352"
$ws.Range("E46").Value = "  +2.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.07039"
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.840"
$ws.Range("E48").Value = "  +1.96%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "112.69"
$ws.Range("E49").Value = "  +1.09%  "
$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.531"
$ws.Range("E50").Value = "  +4.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.068"
$ws.Range("E51").Value = "  -5.44%  "
